$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the ID column to remain text (preserves leading zeros, matches original schema)
$ws.Range("D2:D31").NumberFormat = "@"

$ws.Range("A2").Value = "Examtaker67234"
$ws.Range("B2").Value = "Automation67234"
$ws.Range("C2").Value = "examtakerautomation67234@gmail.com"
$ws.Range("D2").Value = "67234"

$ws.Range("A3").Value = "Examtaker67780"
$ws.Range("B3").Value = "Automation67780"
$ws.Range("C3").Value = "examtakerautomation67780@gmail.com"
$ws.Range("D3").Value = "67780"

$ws.Range("A4").Value = "Examtaker99976"
$ws.Range("B4").Value = "Automation99976"
$ws.Range("C4").Value = "examtakerautomation99976@gmail.com"
$ws.Range("D4").Value = "99976"

$ws.Range("A5").Value = "Examtaker96782"
$ws.Range("B5").Value = "Automation96782"
$ws.Range("C5").Value = "examtakerautomation96782@gmail.com"
$ws.Range("D5").Value = "96782"

$ws.Range("A6").Value = "Examtaker23971"
$ws.Range("B6").Value = "Automation23971"
$ws.Range("C6").Value = "examtakerautomation23971@gmail.com"
$ws.Range("D6").Value = "23971"

$ws.Range("A7").Value = "Examtaker44871"
$ws.Range("B7").Value = "Automation44871"
$ws.Range("C7").Value = "examtakerautomation44871@gmail.com"
$ws.Range("D7").Value = "44871"

$ws.Range("A8").Value = "Examtaker91744"
$ws.Range("B8").Value = "Automation91744"
$ws.Range("C8").Value = "examtakerautomation91744@gmail.com"
$ws.Range("D8").Value = "91744"

$ws.Range("A9").Value = "Examtaker83268"
$ws.Range("B9").Value = "Automation83268"
$ws.Range("C9").Value = "examtakerautomation83268@gmail.com"
$ws.Range("D9").Value = "83268"

$ws.Range("A10").Value = "Examtaker19536"
$ws.Range("B10").Value = "Automation19536"
$ws.Range("C10").Value = "examtakerautomation19536@gmail.com"
$ws.Range("D10").Value = "19536"

$ws.Range("A11").Value = "Examtaker45715"
$ws.Range("B11").Value = "Automation45715"
$ws.Range("C11").Value = "examtakerautomation45715@gmail.com"
$ws.Range("D11").Value = "45715"

$ws.Range("A12").Value = "Examtaker01631"
$ws.Range("B12").Value = "Automation01631"
$ws.Range("C12").Value = "examtakerautomation01631@gmail.com"
$ws.Range("D12").Value = "01631"

$ws.Range("A13").Value = "Examtaker91177"
$ws.Range("B13").Value = "Automation91177"
$ws.Range("C13").Value = "examtakerautomation91177@gmail.com"
$ws.Range("D13").Value = "91177"

$ws.Range("A14").Value = "Examtaker19514"
$ws.Range("B14").Value = "Automation19514"
$ws.Range("C14").Value = "examtakerautomation19514@gmail.com"
$ws.Range("D14").Value = "19514"

$ws.Range("A15").Value = "Examtaker25968"
$ws.Range("B15").Value = "Automation25968"
$ws.Range("C15").Value = "examtakerautomation25968@gmail.com"
$ws.Range("D15").Value = "25968"

$ws.Range("A16").Value = "Examtaker66463"
$ws.Range("B16").Value = "Automation66463"
$ws.Range("C16").Value = "examtakerautomation66463@gmail.com"
$ws.Range("D16").Value = "66463"

$ws.Range("A17").Value = "Examtaker49863"
$ws.Range("B17").Value = "Automation49863"
$ws.Range("C17").Value = "examtakerautomation49863@gmail.com"
$ws.Range("D17").Value = "49863"

$ws.Range("A18").Value = "Examtaker62120"
$ws.Range("B18").Value = "Automation62120"
$ws.Range("C18").Value = "examtakerautomation62120@gmail.com"
$ws.Range("D18").Value = "62120"

$ws.Range("A19").Value = "Examtaker89069"
$ws.Range("B19").Value = "Automation89069"
$ws.Range("C19").Value = "examtakerautomation89069@gmail.com"
$ws.Range("D19").Value = "89069"

$ws.Range("A20").Value = "Examtaker24311"
$ws.Range("B20").Value = "Automation24311"
$ws.Range("C20").Value = "examtakerautomation24311@gmail.com"
$ws.Range("D20").Value = "24311"

$ws.Range("A21").Value = "Examtaker57572"
$ws.Range("B21").Value = "Automation57572"
$ws.Range("C21").Value = "examtakerautomation57572@gmail.com"
$ws.Range("D21").Value = "57572"

$ws.Range("A22").Value = "Examtaker47847"
$ws.Range("B22").Value = "Automation47847"
$ws.Range("C22").Value = "examtakerautomation47847@gmail.com"
$ws.Range("D22").Value = "47847"

$ws.Range("A23").Value = "Examtaker38505"
$ws.Range("B23").Value = "Automation38505"
$ws.Range("C23").Value = "examtakerautomation38505@gmail.com"
$ws.Range("D23").Value = "38505"

$ws.Range("A24").Value = "Examtaker27293"
$ws.Range("B24").Value = "Automation27293"
$ws.Range("C24").Value = "examtakerautomation27293@gmail.com"
$ws.Range("D24").Value = "27293"

$ws.Range("A25").Value = "Examtaker59284"
$ws.Range("B25").Value = "Automation59284"
$ws.Range("C25").Value = "examtakerautomation59284@gmail.com"
$ws.Range("D25").Value = "59284"

$ws.Range("A26").Value = "Examtaker12605"
$ws.Range("B26").Value = "Automation12605"
$ws.Range("C26").Value = "examtakerautomation12605@gmail.com"
$ws.Range("D26").Value = "12605"

$ws.Range("A27").Value = "Examtaker43609"
$ws.Range("B27").Value = "Automation43609"
$ws.Range("C27").Value = "examtakerautomation43609@gmail.com"
$ws.Range("D27").Value = "43609"

$ws.Range("A28").Value = "Examtaker30736"
$ws.Range("B28").Value = "Automation30736"
$ws.Range("C28").Value = "examtakerautomation30736@gmail.com"
$ws.Range("D28").Value = "30736"

$ws.Range("A29").Value = "Examtaker71699"
$ws.Range("B29").Value = "Automation71699"
$ws.Range("C29").Value = "examtakerautomation71699@gmail.com"
$ws.Range("D29").Value = "71699"

$ws.Range("A30").Value = "Examtaker20065"
$ws.Range("B30").Value = "Automation20065"
$ws.Range("C30").Value = "examtakerautomation20065@gmail.com"
$ws.Range("D30").Value = "20065"

$ws.Range("A31").Value = "Examtaker28138"
$ws.Range("B31").Value = "Automation28138"
$ws.Range("C31").Value = "examtakerautomation28138@gmail.com"
$ws.Range("D31").Value = "28138"
